$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Sprint 5 table: update actuals (column F) for a few tasks ---
$ws.Range("F14").Value2 = 35
$ws.Range("F17").Value2 = 35
$ws.Range("F18").Value2 = 3
$ws.Range("F19").Value2 = 3

# --- Sprint 6 (second) table: insert a new planning row for
#     "Container opdracht verbeteren" right after the "Unit tests" row (row 27) ---
$ws.Rows.Item(28).Insert()

# Copy row 27's formatting down onto the freshly inserted row 28 so the new
# row matches the surrounding table styling (borders, fonts, alignment...).
$ws.Range("A27:I27").Copy()
$ws.Range("A28:I28").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Fill in the new row's content.
$ws.Range("A28").Value2 = ""
$ws.Range("B28").Value2 = "Sprint 5"
$ws.Range("C28").Value2 = "Container opdracht verbeteren"
$ws.Range("D28").Value2 = 1
$ws.Range("E28").Value2 = 25
$ws.Range("F28").Value2 = 0
$ws.Range("G28").Formula = "=E28-F28"
$ws.Range("H28").Value2 = ""
$ws.Range("I28").Value2 = ""

# The conditional-formatting blocks that covered the old row range (25:30)
# need to grow by one row too, now that row 28 pushed everything down.
$fcG = $ws.Range("G25:G30").FormatConditions.Item(1)
$fcG.ModifyAppliesToRange($ws.Range("G25:G31"))
$fcD = $ws.Range("D25:D30").FormatConditions.Item(1)
$fcD.ModifyAppliesToRange($ws.Range("D25:D31"))

# --- move the view / selection as recorded in the source file ---
$ws.Range("F19").Select()
